$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '利欧股份'
$ws.Range("C2").Value = '风语筑'
$ws.Range("B3").Value = '协鑫集成'
$ws.Range("C3").Value = '大位科技'
$ws.Range("A4").Value = '风语筑'
$ws.Range("B4").Value = '华胜天成'
$ws.Range("C4").Value = '天奇股份'
$ws.Range("A5").Value = '协鑫集成'
$ws.Range("B5").Value = '湖南白银'
$ws.Range("C5").Value = '协鑫集成'
$ws.Range("A6").Value = '利欧股份'
$ws.Range("C6").Value = '华胜天成'
$ws.Range("A7").Value = '光线传媒'
$ws.Range("B7").Value = '白银有色'
$ws.Range("C7").Value = '光线传媒'
$ws.Range("A8").Value = '国际复材'
$ws.Range("B8").Value = '光线传媒'
$ws.Range("C8").Value = '汉缆股份'
$ws.Range("A9").Value = '汉缆股份'
$ws.Range("B9").Value = '风语筑'
$ws.Range("C9").Value = '利欧股份'
$ws.Range("A10").Value = '大位科技'
$ws.Range("B10").Value = '汉缆股份'
$ws.Range("C10").Value = '中国中免'
$ws.Range("A11").Value = '浙江世宝'
$ws.Range("B11").Value = '浙江世宝'
$ws.Range("C11").Value = '亚星锚链'
$ws.Range("A12").Value = '白银有色'
$ws.Range("B12").Value = '三花智控'
$ws.Range("C12").Value = '博纳影业'
$ws.Range("A13").Value = '湖南白银'
$ws.Range("B13").Value = '华工科技'
$ws.Range("C13").Value = '华工科技'
$ws.Range("A14").Value = '长飞光纤'
$ws.Range("B14").Value = '大位科技'
$ws.Range("C14").Value = '长飞光纤'
$ws.Range("B15").Value = '国际复材'
$ws.Range("C15").Value = '红宝丽'
$ws.Range("A16").Value = '洲际油气'
$ws.Range("B16").Value = '洲际油气'
$ws.Range("C16").Value = '洲际油气'
$ws.Range("A17").Value = '润泽科技'
$ws.Range("B17").Value = '保变电气'
$ws.Range("C17").Value = '巨力索具'
$ws.Range("A18").Value = '中国中免'
$ws.Range("B18").Value = '特变电工'
$ws.Range("C18").Value = '白银有色'
$ws.Range("A19").Value = '红 宝 丽'
$ws.Range("B19").Value = '深科技'
$ws.Range("C19").Value = '再升科技'
$ws.Range("A20").Value = '风华高科'
$ws.Range("B20").Value = '红 宝 丽'
$ws.Range("C20").Value = '掌阅科技'
$ws.Range("A21").Value = '亚星锚链'
$ws.Range("B21").Value = '中国中铁'
$ws.Range("C21").Value = '湖南白银'
